# Auto-generated edit script
$d = $word.ActiveDocument

function New-ListParagraph($afterRange, $level) {
    $afterRange.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Range.ListFormat.ListLevelNumber = $level
    return $p
}

function Append-Run($para, $text) {
    $rng = $d.Range($para.Range.End - 1, $para.Range.End - 1)
    $rng.InsertAfter($text)
}

# 1. "6/20/2023 - Dan" (ilvl 0)
$last = $d.Paragraphs.Last
$p1 = New-ListParagraph $last.Range 1
$p1.Range.InsertAfter("6/20/2023 – Dan")

# 2. Corrected typo paragraph (ilvl 1) - 5 runs
$p2 = New-ListParagraph $d.Paragraphs.Last.Range 2
$p2.Range.InsertAfter("Corrected typo in ui.R of validation app that mis-represented acceptable Wr values as ")
Append-Run $p2 "<60 or >120"
Append-Run $p2 " and corrected to "
Append-Run $p2 "<50 or >150"
Append-Run $p2 ", which is what the app is actually testing."

# 3. "Date-Dan" (ilvl 0)
$p3 = New-ListParagraph $d.Paragraphs.Last.Range 1
$p3.Range.InsertAfter("Date-Dan")

# 4. Force only paragraph (ilvl 1)
$p4 = New-ListParagraph $d.Paragraphs.Last.Range 2
$p4.Range.InsertAfter("Force only “Validated” or “validated” for Verified.TL and Verified.Wr columns.  Added text to right column indicating what to type in these fields as a reminder.")

# 5. Changed download paragraph (ilvl 1) - 3 runs
$p5 = New-ListParagraph $d.Paragraphs.Last.Range 2
$p5.Range.InsertAfter("Changed download of verified data to add lake year and gear so the format ")
Append-Run $p5 "of file name "
Append-Run $p5 "is now: lake_year_gear.code_Sample_verified_date.csv"

# 6. Added clarifying text paragraph (ilvl 1)
$p6 = New-ListParagraph $d.Paragraphs.Last.Range 2
$p6.Range.InsertAfter("Added clarifying text about needing period in any blank cell in right hand column text related to this.")

# 7. Require TL_mm paragraph (ilvl 1)
$p7 = New-ListParagraph $d.Paragraphs.Last.Range 2
$p7.Range.InsertAfter("Require TL_mm (SSP data) and TLmm (Age data) to be integers (S-central region measures in inches and converts to mm...this will require them to round to whole number to pass validation).  This was important as it changes the data type and creates problems when merging into he main databases if these columns are not integer data types.")

Write-Output "Done. Paragraph count:"
Write-Output $d.Paragraphs.Count
